# Updated symbol list on Tue Feb  7 14:23:01 UTC 2023 with GitHub Actions
# Applies the refreshed Price (D) / Volume(1h) (E) values to the cryptos sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# row => @{ D = "new price"; E = "new volume%" }
$updates = @{
    2  = @{ D = "329.09";     E = "0.50%" }
    3  = @{ D = "44.27";      E = "0.52%" }
    4  = @{ D = "5.525";      E = "-0.69%" }
    5  = @{ D = "0.08075";    E = "0.44%" }
    6  = @{ D = "2.058";      E = "8.55%" }
    7  = @{ D = "0.9536";     E = "1.00%" }
    8  = @{ D = "0.1146";     E = "-1.46%" }
    9  = @{ D = "0.1899";     E = "3.19%" }
    10 = @{ D = "10.10";      E = "1.54%" }
    11 = @{ D = "0.09884";    E = "2.41%" }
    12 = @{ D = "0.04861";    E = "10.72%" }
    13 = @{ E = "-0.54%" }
    14 = @{ D = "0.001265";   E = "-1.87%" }
    15 = @{ D = "0.04088";    E = "-2.70%" }
    16 = @{ D = "0.006079";   E = "1.77%" }
    17 = @{ E = "-0.77%" }
    18 = @{ D = "4.409";      E = "2.89%" }
    19 = @{ D = "2.586";      E = "0.73%" }
    21 = @{ E = "0.15%" }
    22 = @{ D = "0.2579";     E = "2.72%" }
    23 = @{ E = "4.20%" }
    24 = @{ E = "1.33%" }
    25 = @{ D = "0.0001250";  E = "-0.88%" }
    26 = @{ D = "0.0003741";  E = "-6.47%" }
    38 = @{ D = "0.02592";    E = "-1.77%" }
    39 = @{ D = "0.05735";    E = "4.70%" }
    40 = @{ D = "0.007582";   E = "0.19%" }
    41 = @{ E = "0.59%" }
    42 = @{ D = "0.007359";   E = "-8.21%" }
    43 = @{ D = "0.002008";   E = "-0.12%" }
    44 = @{ D = "0.009057";   E = "2.54%" }
    45 = @{ D = "0.00007014"; E = "1.17%" }
    46 = @{ E = "-0.10%" }
    47 = @{ D = "0.0005801";  E = "-0.18%" }
    48 = @{ D = "0.003499";   E = "53.72%" }
    49 = @{ E = "-1.63%" }
    50 = @{ E = "-0.10%" }
    51 = @{ E = "-0.10%" }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    if ($rowData.ContainsKey("D")) {
        Set-TextValue $ws.Range("D$row") $rowData["D"]
    }
    if ($rowData.ContainsKey("E")) {
        Set-TextValue $ws.Range("E$row") $rowData["E"]
    }
}
